# Replace the division-problem text in the worksheet table.
# The table has 20 rows; only every 4th row (1, 5, 9, 13, 17) holds the
# five "NN÷N=" problems per row, the rows in between are blank spacer rows.
# Addressing cells by (row, column) avoids ambiguity from duplicate values
# (e.g. "60÷6=" appears twice in the original document).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "27÷4="
$t.Cell(1, 2).Range.Text = "53÷4="
$t.Cell(1, 3).Range.Text = "33÷2="
$t.Cell(1, 4).Range.Text = "31÷9="
$t.Cell(1, 5).Range.Text = "45÷6="

$t.Cell(5, 1).Range.Text = "79÷4="
$t.Cell(5, 2).Range.Text = "73÷7="
$t.Cell(5, 3).Range.Text = "42÷7="
$t.Cell(5, 4).Range.Text = "31÷4="
$t.Cell(5, 5).Range.Text = "64÷8="

$t.Cell(9, 1).Range.Text = "68÷3="
$t.Cell(9, 2).Range.Text = "57÷2="
$t.Cell(9, 3).Range.Text = "47÷8="
$t.Cell(9, 4).Range.Text = "46÷7="
$t.Cell(9, 5).Range.Text = "72÷7="

$t.Cell(13, 1).Range.Text = "29÷9="
$t.Cell(13, 2).Range.Text = "34÷5="
$t.Cell(13, 3).Range.Text = "80÷4="
$t.Cell(13, 4).Range.Text = "74÷7="
$t.Cell(13, 5).Range.Text = "75÷3="

$t.Cell(17, 1).Range.Text = "32÷4="
$t.Cell(17, 2).Range.Text = "35÷5="
$t.Cell(17, 3).Range.Text = "83÷3="
$t.Cell(17, 4).Range.Text = "19÷7="
$t.Cell(17, 5).Range.Text = "37÷4="
